# Monthly rollover update.
#
# "VENTAS POR GRUPO"  -> the current-month (by product-category) figures are
#                        reset to 0 now that the month has closed, and the
#                        "<n> de 57" tally row (row 59) is refreshed for every
#                        column whose figure was just zeroed.
#
# "VENTA MENSUAL"     -> the rolling 4-month window (columns C:F) shifts one
#                        month to the left (oldest month "mayo" drops off,
#                        a brand-new "septiembre" column appears at the far
#                        right seeded with zeros), and the column widths
#                        shift along with the data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Cells whose current-month value goes back to 0.
$zeroedCells = @(
    "H2", "L2", "M2",
    "O7",
    "E11", "G11", "M11", "N11",
    "L13", "M13",
    "D22",
    "L28", "M28", "O28",
    "M32",
    "H33",
    "D45",
    "D46", "M46",
    "M50",
    "O51",
    "M56",
    "M57",
    "M58"
)

foreach ($cellRef in $zeroedCells) {
    $wsGrupo.Range($cellRef).Value = 0
}

# The "<count> de 57" summary in row 59 is recomputed only for the columns
# that just lost their only (or last) non-zero contributor.
$summaryResetCells = @("D59", "E59", "G59", "H59", "L59", "M59", "N59", "O59")

foreach ($cellRef in $summaryResetCells) {
    $wsGrupo.Range($cellRef).Value = "0 de 57"
}

# ---------------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Shift the header months left by one and introduce "septiembre".
$wsMensual.Range("C1").Value = $wsMensual.Range("D1").Value2
$wsMensual.Range("D1").Value = $wsMensual.Range("E1").Value2
$wsMensual.Range("E1").Value = $wsMensual.Range("F1").Value2
$wsMensual.Range("F1").Value = "septiembre"

# Shift every data row's C:F figures left by one column; the new F column
# starts the month at 0.
$lastRow = 59
for ($row = 2; $row -le $lastRow; $row++) {
    $oldC = $wsMensual.Cells.Item($row, 3).Value2
    $oldD = $wsMensual.Cells.Item($row, 4).Value2
    $oldE = $wsMensual.Cells.Item($row, 5).Value2
    $oldF = $wsMensual.Cells.Item($row, 6).Value2

    $wsMensual.Cells.Item($row, 3).Value = $oldD
    $wsMensual.Cells.Item($row, 4).Value = $oldE
    $wsMensual.Cells.Item($row, 5).Value = $oldF
    $wsMensual.Cells.Item($row, 6).Value = 0
}

# The per-column widths travel with the data they describe; the new,
# wider "septiembre" column is sized to fit its longer header text.
# (15.2 is within the pixel-rounding plateau that serialises to a raw
# OOXML column width of exactly 16.)
$wsMensual.Columns.Item(4).ColumnWidth = $wsMensual.Columns.Item(5).ColumnWidth
$wsMensual.Columns.Item(5).ColumnWidth = $wsMensual.Columns.Item(6).ColumnWidth
$wsMensual.Columns.Item(6).ColumnWidth = 15.2
